# Add additional parsed OCR receipt rows (and normalise the numeric
# columns on the existing row) so the sheet can hold multiple .pdf
# receipts instead of just one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: Bill Number / Total Amount were imported as text; they are
# genuinely numeric, so store them as real numbers. ---
$ws.Range("B2").Value = 258
$ws.Range("D2").Value = 13.32

# --- Row 3: duplicate of the first receipt (re-parsed from the .pdf) ---
$ws.Range("A3").Value = "THE BREAKFAST CLUB"
$ws.Range("B3").Value = 258
$ws.Range("C3").Value = ""
$ws.Range("D3").Value = 13.32
$ws.Range("E3").Value = "5000 NORTH ALAMAR AVE"

# --- Row 4: Montana Restaurant receipt ---
$ws.Range("A4").Value = "Montana Restaurant."
$ws.Range("B4").Value = ""
$ws.Range("C4").NumberFormat = "@"
$ws.Range("C4").Value = "10/07/2020"
$ws.Range("D4").Value = 36.98
$ws.Range("E4").Value = "6542 MAGNOLIA LAKE COURT"

# --- Row 5: Harbor Lane Cafe receipt ---
$ws.Range("A5").Value = "HARBOR LANE CAFE"
$ws.Range("B5").Value = ""
$ws.Range("C5").NumberFormat = "@"
$ws.Range("C5").Value = "11/20/2019"
$ws.Range("D5").Value = 31.39
$ws.Range("E5").Value = "3941 GREEN OAKS BLVD"

# --- Row 6: same Harbor Lane Cafe receipt, re-OCR'd with the amount
# still as text (not yet normalised to a number) ---
$ws.Range("A6").Value = "HARBOR LANE CAFE"
$ws.Range("B6").Value = ""
$ws.Range("C6").NumberFormat = "@"
$ws.Range("C6").Value = "11/20/2019"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "31.39"
$ws.Range("E6").Value = "3941 GREEN OAKS BLVD"
